# Applies the "Added plots for number of facilities and cumulative investment"
# edit: adds a new "Export cable" worksheet (mirroring the Monopile/Blade
# layout), appends four rows to the "Avg Demand Scenario" summary table, tweaks
# two values on the "Blade" sheet, and switches the active tab to "Blade".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "Export cable" worksheet right after "Blade" (so the tab
#    order becomes Avg Demand Scenario, Monopile, Blade, Export cable).
# ---------------------------------------------------------------------------
$blade = $wb.Worksheets.Item("Blade")
$cable = $wb.Worksheets.Add($null, $blade)
$cable.Name = "Export cable"

# Column layout (set first so AutoFit below has content to measure against).
$cable.Range("A1").Value = "Factory specifications"
$cable.Range("A1").Font.Bold = $true
$cable.Range("B1").Value = "Value"
$cable.Range("C1").Value = "Unit"

$cable.Range("A2").Value = "Annual throughput"
$cable.Range("B2").Value = 200
$cable.Range("C2").Value = "km/year"

$cable.Range("A3").Value = "Investment cost"
$cable.Range("B3").Value = 200
$cable.Range("C3").Value = "USD"

$cable.Range("A4").Value = "Lead time"
$cable.Range("B4").Value = 6
$cable.Range("C4").Value = "years"

$cable.Range("A5").Value = "Workforce skills assessment"
$cable.Range("A5").Font.Bold = $true

$cable.Range("A6").Value = "Total Direct Jobs"
$cable.Range("B6").Value = 200
$cable.Range("C6").Value = "FTE"

$cable.Range("A7").Value = "Regional Professionals"
$cable.Range("B7").Value = 10
$cable.Range("C7").Value = "%"

$cable.Range("A8").Value = "Factory-Level Management"
$cable.Range("B8").Value = 20
$cable.Range("C8").Value = "%"

$cable.Range("A9").Value = "Design and Engineering"
$cable.Range("B9").Value = 15
$cable.Range("C9").Value = "%"

$cable.Range("A10").Value = "Quality and Safety"
$cable.Range("B10").Value = 20
$cable.Range("C10").Value = "%"

$cable.Range("A11").Value = "Factory-Level Worker"
$cable.Range("B11").Value = 30
$cable.Range("C11").Value = "%"

$cable.Range("A12").Value = "Facilities Maintenance"
$cable.Range("B12").Value = 5
$cable.Range("C12").Value = "%"

$cable.Columns.AutoFit()
$cable.Range("A16").Select()

# ---------------------------------------------------------------------------
# 2) Extend the "Avg Demand Scenario" table with the new facility rows.
# ---------------------------------------------------------------------------
$avg = $wb.Worksheets.Item("Avg Demand Scenario")
$avg.Activate()

$avg.Range("A4").Value = "Monopile 1"
$avg.Range("B4").Value = 2026
$avg.Range("C4").Value = "MA"

$avg.Range("A5").Value = "Export cable 1"
$avg.Range("B5").Value = 2028

$avg.Range("A6").Value = "Export cable 2"
$avg.Range("B6").Value = 2029

$avg.Range("A7").Value = "Export cable 3"
$avg.Range("B7").Value = 2030

$avg.Range("B7").Select()

# ---------------------------------------------------------------------------
# 3) Update the "Blade" sheet's figures and make it the active tab.
# ---------------------------------------------------------------------------
$blade.Range("B2").Value = 200
$blade.Range("B4").Value = 4

$blade.Activate()
$blade.Range("E7").Select()
